$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.746.97'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.907.22'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('E3').ClearFormats()

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9986'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.54'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9985'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E6').ClearFormats()

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5222'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +7.32%  '
$ws.Range('E7').ClearFormats()

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3782'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07243'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.34'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.71%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9014'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07639'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.69%  '
$ws.Range('E12').ClearFormats()

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.861.65'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('E13').ClearFormats()

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.454'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.12'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9992'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008704'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9991'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('E18').ClearFormats()

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '27.766.12'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.51'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('E20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.145'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.46%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.138.46'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.84'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('E23').ClearFormats()

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.587'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('E24').ClearFormats()

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.05'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('E25').ClearFormats()

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.880'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.171'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.31'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.52'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.862'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08981'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('E31').ClearFormats()

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.863'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.68%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.176'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.232'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7737'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.13%  '
$ws.Range('E35').ClearFormats()

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.638'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.57%  '
$ws.Range('E36').ClearFormats()

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02084'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.03%  '
$ws.Range('E37').ClearFormats()

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.070'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.81%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.094'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('E39').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5512'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('E40').ClearFormats()

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.05288'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.668'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.50%  '
$ws.Range('E42').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '114.56'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.16%  '
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.519'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1510'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4810'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('E46').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.44'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.79%  '
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9981'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.617'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '66.78'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05986'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.14%  '
$ws.Range('E51').ClearFormats()
